# Swap the data between row 2 and row 3 for columns A, Q, R, S, AW, AX
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "Q", "R", "S", "AW", "AX")

foreach ($col in $cols) {
    $range2 = $ws.Range($col + "2")
    $range3 = $ws.Range($col + "3")

    $val2 = $range2.Value2
    $val3 = $range3.Value2

    $range2.Value = $val3
    $range3.Value = $val2
}
